$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 181, shifting existing rows 181-268 down to 182-269.
$ws.Rows.Item(181).Insert()

# Populate the newly inserted row 181 with the new record.
$ws.Cells.Item(181, 1).Value = 10
$ws.Cells.Item(181, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(181, 3).Value = "La Araucanía"
$ws.Cells.Item(181, 4).Value = 44529
$ws.Cells.Item(181, 5).Value = 9
$ws.Cells.Item(181, 6).Value = 100112040
$ws.Cells.Item(181, 7).Value = "Cilantro"
$ws.Cells.Item(181, 8).Value = "Sin especificar"
$ws.Cells.Item(181, 9).Value = "Primera"
$ws.Cells.Item(181, 10).Value = 90
$ws.Cells.Item(181, 11).Value = 4500
$ws.Cells.Item(181, 12).Value = 5000
$ws.Cells.Item(181, 13).Value = 4694
$ws.Cells.Item(181, 14).Value = "$/docena de atados (2 kilos)"
$ws.Cells.Item(181, 15).Value = "Provincia de Cautín"
$ws.Cells.Item(181, 16).Value = 2347
$ws.Cells.Item(181, 17).Value = 2
$ws.Cells.Item(181, 18).Value = "Hortaliza"
